{"js": "const body = context.document.body;\nconst searchResults = body.search(\"Dott. Beniamino Gioli\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"{{ nome_cognome_direttore}}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n[void]$rng.Find.Execute(\"Dott. Beniamino Gioli\", $false, $false, $false, $false, $false, $true, 0, $false, \"{{ nome_cognome_direttore}}\", 2)\n"}
